$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.320.28'
$ws.Range('E2').Value = '  -4.92%  '
$ws.Range('D3').Value = '3.244.32'
$ws.Range('E3').Value = '  -6.29%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '552.75'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -4.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '178.62'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -6.38%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.583'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -3.30%  '
$ws.Range('D9').Value = '3.233.26'
$ws.Range('E9').Value = '  -6.29%  '
$ws.Range('E10').Value = '  -8.96%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.582'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -5.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.08'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -8.15%  '
$ws.Range('E13').Value = '  -7.42%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '629.19'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.50'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -6.21%  '
$ws.Range('D16').Value = '3.753.77'
$ws.Range('E16').Value = '  -6.20%  '
$ws.Range('D17').Value = '65.110.88'
$ws.Range('E17').Value = '  -4.91%  '
$ws.Range('E18').Value = '  -2.20%  '
$ws.Range('E19').Value = '  -3.35%  '
$ws.Range('D20').Value = '3.243.02'
$ws.Range('E20').Value = '  -6.23%  '
$ws.Range('E21').Value = '  -8.45%  '
$ws.Range('E22').Value = '  -4.66%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.58'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '105.65'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +6.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.93'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -7.52%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.94'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -8.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.65'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -6.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.47'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -3.15%  '
$ws.Range('E29').Value = '  -6.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.15'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -6.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.95'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -4.94%  '
$ws.Range('E32').Value = '  -6.93%  '
$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '10.99'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -5.01%  '
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '541.64'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +8.46%  '
$ws.Range('E35').Value = '  -3.59%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '56.83'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -6.87%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').Value = '3.573.49'
$ws.Range('E38').Value = '  -2.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.67'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +8.17%  '
$ws.Range('E40').Value = '  -1.51%  '
$ws.Range('E41').Value = '  -5.97%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.128'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -3.57%  '
$ws.Range('B43').Value = 'PEPE'
$ws.Range('C43').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D43').Value = '0.0₃0701'
$ws.Range('E43').Value = '  -9.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '31.82'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -7.78%  '
$ws.Range('E45').Value = '  -9.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.28'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.76%  '
$ws.Range('E47').Value = '  -6.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.128'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.85%  '
$ws.Range('E49').Value = '  -7.94%  '
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.23'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.84%  '
